# "Tested home page on Windows"
# Fill in the "Windows OS" result column (G) on the HomePage sheet with
# "Works" for each executed test case (all test cases except test case #4,
# which was left unmarked).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomePage")

$rows = @(2,3,4,6,7,8,9,10,11,12,13,14,15)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = "Works"
}
